# Refresh the "Price" (D) and "Volume(1h)" (E) text columns of the
# cryptos worksheet with the latest feed values, including the
# PancakeSwap / Filecoin rows trading ranking positions (their Coin,
# Link, Price and Volume cells swap).
#
# Every cell in columns B-E is stored as literal text (e.g. "1.00",
# "34.118.08", "  -1.73%  "), never a number. Writing a numeric-looking
# literal straight to Range.Value lets Excel "helpfully" reinterpret it
# as a real number -- dropping significant trailing zeros (e.g. "1.00"
# -> 1) and upgrading the cell to a new quote-prefixed-text style, which
# would not match the source workbook. To avoid that, numeric-looking
# replacements are staged in a scratch cell as a ="literal" text formula
# and brought into the destination with PasteSpecial(xlPasteValues) --
# which lands an exact text value while leaving the destination cell's
# original (default) style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163
$scratch = $ws.Range('ZZ1')

function Set-TextCell([string]$cellRef, [string]$text) {
    $target = $ws.Range($cellRef)
    $looksNumeric = ($text.Trim() -match '^[+-]?\d+(\.\d+)?$')
    if ($looksNumeric) {
        $escaped = $text.Replace('"', '""')
        $scratch.Formula = '="' + $escaped + '"'
        $scratch.Copy()
        $target.PasteSpecial($xlPasteValues)
    } else {
        $target.Value = $text
    }
}


Set-TextCell 'D2' '34.118.08'
Set-TextCell 'E2' '  -1.73%  '
Set-TextCell 'D3' '1.798.09'
Set-TextCell 'E3' '  +0.44%  '
Set-TextCell 'D4' '1.00'
Set-TextCell 'E4' '  -0.03%  '
Set-TextCell 'D5' '222.84'
Set-TextCell 'E5' '  -0.15%  '
Set-TextCell 'E6' '  -0.71%  '
Set-TextCell 'D7' '0.999'
Set-TextCell 'E7' '  -0.01%  '
Set-TextCell 'D8' '32.33'
Set-TextCell 'E8' '  -0.61%  '
Set-TextCell 'E9' '  +1.64%  '
Set-TextCell 'D10' '0.0716'
Set-TextCell 'E10' '  +3.97%  '
Set-TextCell 'E11' '  -1.36%  '
Set-TextCell 'D12' '2.055.54'
Set-TextCell 'E12' '  +0.38%  '
Set-TextCell 'D13' '1.793.51'
Set-TextCell 'E13' '  +0.22%  '
Set-TextCell 'D14' '10.72'
Set-TextCell 'E14' '  -2.41%  '
Set-TextCell 'E15' '  +0.18%  '
Set-TextCell 'D16' '34.134.81'
Set-TextCell 'E16' '  -1.69%  '
Set-TextCell 'D17' '4.21'
Set-TextCell 'E17' '  -1.79%  '
Set-TextCell 'D18' '68.20'
Set-TextCell 'E18' '  -0.46%  '
Set-TextCell 'D19' '246.71'
Set-TextCell 'E19' '  -2.64%  '
Set-TextCell 'E20' '  -0.16%  '
Set-TextCell 'D21' '1.00'
Set-TextCell 'E21' '  +0.02%  '
Set-TextCell 'D22' '10.87'
Set-TextCell 'E22' '  +3.94%  '
Set-TextCell 'E23' '  -1.98%  '
Set-TextCell 'E24' '  -0.39%  '
Set-TextCell 'D25' '159.03'
Set-TextCell 'E25' '  +0.48%  '
Set-TextCell 'D26' '16.56'
Set-TextCell 'E26' '  +1.20%  '
Set-TextCell 'D27' '7.07'
Set-TextCell 'E27' '  +0.36%  '
Set-TextCell 'E28' '  -1.41%  '
Set-TextCell 'D29' '1.00'
Set-TextCell 'E29' '  +0.03%  '
Set-TextCell 'D30' '0.0523'
Set-TextCell 'E30' '  +1.25%  '
Set-TextCell 'B31' 'Filecoin'
Set-TextCell 'C31' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D31' '3.72'
Set-TextCell 'E31' '  -0.73%  '
Set-TextCell 'B32' 'PancakeSwap'
Set-TextCell 'C32' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 'D32' '1.21'
Set-TextCell 'E32' '  +1.31%  '
Set-TextCell 'E33' '  -1.57%  '
Set-TextCell 'E34' '  -1.53%  '
Set-TextCell 'D35' '1.412.90'
Set-TextCell 'E35' '  -1.26%  '
Set-TextCell 'D36' '0.646'
Set-TextCell 'E36' '  +2.18%  '
Set-TextCell 'E37' '  +0.15%  '
Set-TextCell 'E38' '  -1.37%  '
Set-TextCell 'E39' '  +4.86%  '
Set-TextCell 'D40' '80.42'
Set-TextCell 'E40' '  -3.56%  '
Set-TextCell 'E41' '  -2.64%  '
Set-TextCell 'D42' '2.35'
Set-TextCell 'E42' '  -0.26%  '
Set-TextCell 'E43' '  +4.63%  '
Set-TextCell 'E44' '  +0.15%  '
Set-TextCell 'E45' '  -1.10%  '
Set-TextCell 'D46' '107.07'
Set-TextCell 'E46' '  +3.22%  '
Set-TextCell 'D47' '1.954.66'
Set-TextCell 'E47' '  +0.33%  '
Set-TextCell 'E48' '  -2.91%  '
Set-TextCell 'D49' '12.09'
Set-TextCell 'E49' '  +0.82%  '
Set-TextCell 'E51' '  +1.17%  '

$scratch.ClearContents()
$excel.CutCopyMode = $false
